# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 22:22"

# --- Row 4: Estados Unidos (refreshed totals) ---
$ws.Range("B4").Value = 360901
$ws.Range("C4").Value = 24228
$ws.Range("E4").Value = 330902
$ws.Range("G4").Value = 1075
$ws.Range("H4").Value = 10691

# --- Row 7: Alemania (refreshed totals) ---
$ws.Range("B7").Value = 102024
$ws.Range("C7").Value = 1901
$ws.Range("E7").Value = 71629
$ws.Range("G7").Value = 111
$ws.Range("H7").Value = 1695

# --- Rows 18-19: Brasil overtakes Portugal in ranking ---
# Row 18 becomes Brasil with freshly updated figures
$ws.Range("A18").Value = "Brasil"
$ws.Range("B18").Value = 12056
$ws.Range("C18").Value = 802
$ws.Range("D18").Value = 127
$ws.Range("E18").Value = 11376
$ws.Range("F18").Value = 296
$ws.Range("G18").Value = 67
$ws.Range("H18").Value = 553

# Row 19 becomes Portugal, carrying the previous Portugal figures
$ws.Range("A19").Value = "Portugal"
$ws.Range("B19").Value = 11730
$ws.Range("C19").Value = 452
$ws.Range("D19").Value = 140
$ws.Range("E19").Value = 11279
$ws.Range("F19").Value = 270
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = 311

# --- Row 32: Rumania (refreshed totals) ---
$ws.Range("E32").Value = 3475
$ws.Range("G32").Value = 25
$ws.Range("H32").Value = 176

# --- Rows 57-59: Egipto overtakes Ucrania and Croacia in ranking ---
# Row 57 becomes Egipto with freshly updated figures
$ws.Range("A57").Value = "Egipto"
$ws.Range("B57").Value = 1322
$ws.Range("C57").Value = 149
$ws.Range("D57").Value = 259
$ws.Range("E57").Value = 978
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 7
$ws.Range("H57").Value = 85

# Row 58 becomes Ucrania, carrying the previous Ucrania figures
$ws.Range("A58").Value = "Ucrania"
$ws.Range("B58").Value = 1319
$ws.Range("C58").Value = 11
$ws.Range("D58").Value = 28
$ws.Range("E58").Value = 1253
$ws.Range("F58").Value = 16
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 38

# Row 59 becomes Croacia, carrying the previous Croacia figures
$ws.Range("A59").Value = "Croacia"
$ws.Range("B59").Value = 1222
$ws.Range("C59").Value = 40
$ws.Range("D59").Value = 130
$ws.Range("E59").Value = 1076
$ws.Range("F59").Value = 36
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 16

# --- Row 84: Principado de Andorra (refreshed totals) ---
$ws.Range("B84").Value = 525
$ws.Range("C84").Value = 24
$ws.Range("D84").Value = 31
$ws.Range("E84").Value = 473
$ws.Range("G84").Value = 3
$ws.Range("H84").Value = 21

# --- Row 91: Afganistan (refreshed totals) ---
$ws.Range("E91").Value = 338
$ws.Range("G91").Value = 4
$ws.Range("H91").Value = 11

# --- Row 92: Cuba (refreshed totals) ---
$ws.Range("F92").Value = 12
